# edit.ps1 -- apply the tracked changes to consignesProjetJEE PL-SQL.docx
#
# Summary of content-level changes (see commit diff):
#  1. Remove the "_GoBack" bookmark that currently sits in the title
#     paragraph ("Bloc 3: Applications informatiques (JEE)").
#  2. Remove two of the five blank, justified paragraphs that sit right
#     before the "Anne Vandevorst" signature line.
#  3. Remove the stray <w:lastRenderedPageBreak/> marker that precedes
#     "D'utiliser les curseurs, ..." (done by rewriting that run's text).
#  4. Remove the four trailing blank paragraphs after "De documenter les
#     differents scripts d'acces a la base de donnees" and re-create the
#     "_GoBack" bookmark as a zero-length bookmark at the end of that
#     paragraph's text.
#  5. In the header, retag the picture run's East Asian language from
#     fr-BE to ja-JP.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the _GoBack bookmark currently anchored in the title paragraph.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2) Delete two of the blank "jc=both" paragraphs that precede the
#    right-aligned "Anne Vandevorst" paragraph. They are located by
#    scanning backwards from that paragraph.
# ---------------------------------------------------------------------
$anneIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Anne*") {
        $anneIndex = $i
        break
    }
}

if ($anneIndex -eq 0) {
    throw "Could not locate the 'Anne Vandevorst' paragraph"
}

# Remove the two immediately preceding (blank) paragraphs.
$d.Paragraphs.Item($anneIndex - 1).Range.Delete()
$d.Paragraphs.Item($anneIndex - 2).Range.Delete()

# ---------------------------------------------------------------------
# 3) Remove <w:lastRenderedPageBreak/> before "D'utiliser les curseurs...".
#    Re-typing the run's text forces the stale rendering marker to drop.
# ---------------------------------------------------------------------
$curseurText = [char]0x2019
$curseurText = "D" + $curseurText + "utiliser les curseurs, des variables de type record et des tableaux"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("D") -and $p.Range.Text -like "*utiliser les curseurs*") {
        $p.Range.Text = $curseurText
        break
    }
}

# ---------------------------------------------------------------------
# 4) Delete the four trailing blank paragraphs after "De documenter..."
#    and move the _GoBack bookmark (zero-length) to the end of that
#    paragraph's text.
# ---------------------------------------------------------------------
$docIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "De documenter*") {
        $docIndex = $i
        break
    }
}

if ($docIndex -eq 0) {
    throw "Could not locate the 'De documenter...' paragraph"
}

# The four following paragraphs are blank / list placeholders - drop them.
for ($n = 0; $n -lt 4; $n++) {
    $d.Paragraphs.Item($docIndex + 1).Range.Delete()
}

# Re-anchor the _GoBack bookmark at the end of the "De documenter..."
# paragraph's text (before the paragraph mark). A temporary marker
# character is used so the new Range has non-zero length (required for
# Bookmarks.Add to resolve the correct location); once bookmarked, the
# marker text is cleared, leaving a proper zero-length bookmark.
$p = $d.Paragraphs.Item($docIndex)
$r = $p.Range.Duplicate()
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter([char]0x2603)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
$markerBm = $d.Bookmarks.Item("_GoBack")
$markerBm.Range.Text = ""

# ---------------------------------------------------------------------
# 5) Header: retag the picture run's East Asian language fr-BE -> ja-JP.
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
    $shp = $hdr.Range.InlineShapes.Item(1)
    $shp.Range.LanguageIDFarEast = "ja-JP"
}

Write-Output "done"
